$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029527156703935
$ws.Range("D2").Value = 1.047897262884342
$ws.Range("E2").Value = 1.02935018692976
$ws.Range("F2").Value = 1.052629636409509
$ws.Range("I2").Value = 1.04074322205546
$ws.Range("J2").Value = 1.03467343925282
$ws.Range("K2").Value = 1.050658617626625
$ws.Range("L2").Value = 1.032164353328869
$ws.Range("M2").Value = 1.055377839119631
$ws.Range("N2").Value = 1.036142794215261

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03041339205753
$ws.Range("D3").Value = 1.048549689405667
$ws.Range("E3").Value = 1.030101221923502
$ws.Range("F3").Value = 1.053474526229937
$ws.Range("I3").Value = 1.040967879586353
$ws.Range("J3").Value = 1.03520106765096
$ws.Range("K3").Value = 1.051123394108979
$ws.Range("L3").Value = 1.032723669025731
$ws.Range("M3").Value = 1.056035530625677
$ws.Range("N3").Value = 1.036671171906247

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030987354571439
$ws.Range("D4").Value = 1.04897208449352
$ws.Range("E4").Value = 1.030588014556062
$ws.Range("F4").Value = 1.054021970766596
$ws.Range("I4").Value = 1.04111207362751
$ws.Range("J4").Value = 1.035542356000641
$ws.Range("K4").Value = 1.051423680766945
$ws.Range("L4").Value = 1.033085743098407
$ws.Range("M4").Value = 1.05646119811624
$ws.Range("N4").Value = 1.03701294492453

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031228768949245
$ws.Range("D5").Value = 1.049149713059894
$ws.Range("E5").Value = 1.030792857554782
$ws.Range("F5").Value = 1.054252292793535
$ws.Range("I5").Value = 1.04117241086476
$ws.Range("J5").Value = 1.035685803166516
$ws.Range("K5").Value = 1.051549811088887
$ws.Range("L5").Value = 1.033237995985909
$ws.Range("M5").Value = 1.056640170408564
$ws.Range("N5").Value = 1.037156595801834

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031269310515782
$ws.Range("D6").Value = 1.049179540774152
$ws.Range("E6").Value = 1.030827263006107
$ws.Range("F6").Value = 1.05429097516241
$ws.Range("I6").Value = 1.041182525196682
$ws.Range("J6").Value = 1.035709886758095
$ws.Range("K6").Value = 1.051570982425352
$ws.Range("L6").Value = 1.033263562046881
$ws.Range("M6").Value = 1.056670221869626
$ws.Range("N6").Value = 1.037180713594876

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030990579890181
$ws.Range("D7").Value = 1.04897445776643
$ws.Range("E7").Value = 1.03059075091255
$ws.Range("F7").Value = 1.054025047650207
$ws.Range("I7").Value = 1.041112880965428
$ws.Range("J7").Value = 1.035544272869779
$ws.Range("K7").Value = 1.051425366561091
$ws.Range("L7").Value = 1.033087777365292
$ws.Range("M7").Value = 1.056463589469555
$ws.Range("N7").Value = 1.037014864515843

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029826557912735
$ws.Range("D8").Value = 1.048117704496247
$ws.Range("E8").Value = 1.029603831402106
$ws.Range("F8").Value = 1.052915015650669
$ws.Range("I8").Value = 1.040819388817003
$ws.Range("J8").Value = 1.034851778427493
$ws.Range("K8").Value = 1.050815784399343
$ws.Range("L8").Value = 1.032353342947382
$ws.Range("M8").Value = 1.055600087993976
$ws.Range("N8").Value = 1.036321386652018

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027779358726037
$ws.Range("D9").Value = 1.046609846067927
$ws.Range("E9").Value = 1.02787112511815
$ws.Range("F9").Value = 1.050964779796188
$ws.Range("I9").Value = 1.040293258221432
$ws.Range("J9").Value = 1.033630626330165
$ws.Range("K9").Value = 1.049738197224372
$ws.Range("L9").Value = 1.031060450367756
$ws.Range("M9").Value = 1.054079301062544
$ws.Range("N9").Value = 1.035098500378644

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026417296788199
$ws.Range("D10").Value = 1.04560595519979
$ws.Range("E10").Value = 1.026720366639369
$ws.Range("F10").Value = 1.0496686201824
$ws.Range("I10").Value = 1.039936526772612
$ws.Range("J10").Value = 1.032815994440163
$ws.Range("K10").Value = 1.049017580290001
$ws.Range("L10").Value = 1.030199451216849
$ws.Range("M10").Value = 1.053066082850552
$ws.Range("N10").Value = 1.034282711617919

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025828174340049
$ws.Range("D11").Value = 1.045171601029418
$ws.Range("E11").Value = 1.026223134110194
$ws.Range("F11").Value = 1.049108340045995
$ws.Range("I11").Value = 1.039780649989896
$ws.Range("J11").Value = 1.032463137986732
$ws.Range("K11").Value = 1.048705033652389
$ws.Range("L11").Value = 1.02982686461728
$ws.Range("M11").Value = 1.05262751947718
$ws.Range("N11").Value = 1.033929354067851

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025609448359398
$ws.Range("D12").Value = 1.045010314975294
$ws.Range("E12").Value = 1.02603859969004
$ws.Range("F12").Value = 1.048900373974792
$ws.Range("I12").Value = 1.039722539376374
$ws.Range("J12").Value = 1.03233205514863
$ws.Range("K12").Value = 1.048588863879165
$ws.Range("L12").Value = 1.029688505385443
$ws.Range("M12").Value = 1.052464644138104
$ws.Range("N12").Value = 1.033798085077083

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025656361264345
$ws.Range("D13").Value = 1.045044908980804
$ws.Range("E13").Value = 1.026078175688881
$ws.Range("F13").Value = 1.04894497672486
$ws.Range("I13").Value = 1.039735013837893
$ws.Range("J13").Value = 1.032360173579499
$ws.Range("K13").Value = 1.048613786109802
$ws.Range("L13").Value = 1.029718182260939
$ws.Range("M13").Value = 1.052499580221281
$ws.Range("N13").Value = 1.033826243439349

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025810092328274
$ws.Range("D14").Value = 1.04515826799424
$ws.Range("E14").Value = 1.026207877166661
$ws.Range("F14").Value = 1.049091146497851
$ws.Range("I14").Value = 1.039775850849142
$ws.Range("J14").Value = 1.032452302963335
$ws.Range("K14").Value = 1.048695432566292
$ws.Range("L14").Value = 1.029815427058218
$ws.Range("M14").Value = 1.052614055598656
$ws.Range("N14").Value = 1.033918503657477

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025904824442802
$ws.Range("D15").Value = 1.045228119223177
$ws.Range("E15").Value = 1.02628781176154
$ws.Range("F15").Value = 1.049181226042808
$ws.Range("I15").Value = 1.039800983937912
$ws.Range("J15").Value = 1.032509064798502
$ws.Range("K15").Value = 1.048745727606166
$ws.Range("L15").Value = 1.029875347596526
$ws.Range("M15").Value = 1.05268459123167
$ws.Range("N15").Value = 1.03397534610096

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02645640890477
$ws.Range("D16").Value = 1.045634789110624
$ws.Range("E16").Value = 1.026753388657601
$ws.Range("F16").Value = 1.049705824631044
$ws.Range("I16").Value = 1.039946842164925
$ws.Range("J16").Value = 1.0328394100226
$ws.Range("K16").Value = 1.049038312246754
$ws.Range("L16").Value = 1.030224183547472
$ws.Range("M16").Value = 1.053095192505153
$ws.Range("N16").Value = 1.034306160453168

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026802580686808
$ws.Range("D17").Value = 1.045889973911239
$ws.Range("E17").Value = 1.0270457159265
$ws.Range("F17").Value = 1.050035151277513
$ws.Range("I17").Value = 1.040037958291679
$ws.Range("J17").Value = 1.033046596721445
$ws.Range("K17").Value = 1.049221705966429
$ws.Range("L17").Value = 1.030443062012951
$ws.Range("M17").Value = 1.053352797598092
$ws.Range("N17").Value = 1.03451364138089

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027004560364005
$ws.Range("D18").Value = 1.046038851304879
$ws.Range("E18").Value = 1.027216327015157
$ws.Range("F18").Value = 1.050227334973722
$ws.Range("I18").Value = 1.040090968719322
$ws.Range("J18").Value = 1.033167434002151
$ws.Range("K18").Value = 1.049328626667384
$ws.Range("L18").Value = 1.030570752394822
$ws.Range("M18").Value = 1.053503070279212
$ws.Range("N18").Value = 1.034634650264395

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027073440952229
$ws.Range("D19").Value = 1.046089620093518
$ws.Range("E19").Value = 1.027274518161209
$ws.Range("F19").Value = 1.050292880371476
$ws.Range("I19").Value = 1.040109020802801
$ws.Range("J19").Value = 1.033208634456137
$ws.Range("K19").Value = 1.049365075384815
$ws.Range("L19").Value = 1.030614295256874
$ws.Range("M19").Value = 1.053554312044414
$ws.Range("N19").Value = 1.034675909227751

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026765433150005
$ws.Range("D20").Value = 1.045862591631415
$ws.Range("E20").Value = 1.0270143414629
$ws.Range("F20").Value = 1.049999808008588
$ws.Range("I20").Value = 1.040028196462226
$ws.Range("J20").Value = 1.033024368701511
$ws.Range("K20").Value = 1.049202034692724
$ws.Range("L20").Value = 1.030419576112199
$ws.Range("M20").Value = 1.053325157344472
$ws.Range("N20").Value = 1.034491381794618

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025764819581882
$ws.Range("D21").Value = 1.045124885128864
$ws.Range("E21").Value = 1.026169678887805
$ws.Range("F21").Value = 1.049048099072322
$ws.Range("I21").Value = 1.039763831186706
$ws.Range("J21").Value = 1.032425173595476
$ws.Range("K21").Value = 1.048671391808373
$ws.Range("L21").Value = 1.029786789885135
$ws.Range("M21").Value = 1.052580344705702
$ws.Range("N21").Value = 1.033891335762804

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025136275209615
$ws.Range("D22").Value = 1.0446613639532
$ws.Range("E22").Value = 1.025639531815348
$ws.Range("F22").Value = 1.048450572831547
$ws.Range("I22").Value = 1.039596393631674
$ws.Range("J22").Value = 1.032048342437198
$ws.Range("K22").Value = 1.048337316288082
$ws.Range("L22").Value = 1.029389141071249
$ws.Range("M22").Value = 1.052112206255297
$ws.Range("N22").Value = 1.033513969461058

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025469422969199
$ws.Range("D23").Value = 1.044907055877491
$ws.Range("E23").Value = 1.025920484516861
$ws.Range("F23").Value = 1.048767251517928
$ws.Range("I23").Value = 1.039685270894593
$ws.Range("J23").Value = 1.032248116248614
$ws.Range("K23").Value = 1.048514457247175
$ws.Range("L23").Value = 1.029599922016912
$ws.Range("M23").Value = 1.052360360108077
$ws.Range("N23").Value = 1.033714026974201

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026782218326703
$ws.Range("D24").Value = 1.045874964405064
$ws.Range("E24").Value = 1.027028517918959
$ws.Range("F24").Value = 1.050015777823241
$ws.Range("I24").Value = 1.04003260783327
$ws.Range("J24").Value = 1.033034412621699
$ws.Range("K24").Value = 1.04921092344845
$ws.Range("L24").Value = 1.030430188311118
$ws.Range("M24").Value = 1.053337646736102
$ws.Range("N24").Value = 1.034501439978325

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028308131971289
$ws.Range("D25").Value = 1.046999433611614
$ws.Range("E25").Value = 1.028318306330372
$ws.Range("F25").Value = 1.051468265944262
$ws.Range("I25").Value = 1.040430332608112
$ws.Range("J25").Value = 1.033946421814584
$ws.Range("K25").Value = 1.050017177275407
$ws.Range("L25").Value = 1.03139453559044
$ws.Range("M25").Value = 1.054472354992832
$ws.Range("N25").Value = 1.035414744328873
